$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 2): add new columns D/E/F next to existing C "Baseline" ---
$ws.Range("D2").Value = "Waypart"
$ws.Range("E2").Value = "BLISS"
$ws.Range("F2").Value = "Custom"
$ws.Range("D2:F2").Font.Bold = $true

# --- Benchmark rows 3-6: new columns D/E/F each = 1 ---
$ws.Range("D3").Value = 1
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 1

$ws.Range("D4").Value = 1
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 1

$ws.Range("D5").Value = 1
$ws.Range("E5").Value = 1
$ws.Range("F5").Value = 1

$ws.Range("D6").Value = 1
$ws.Range("E6").Value = 1
$ws.Range("F6").Value = 1

# --- Row 8 (CPU # header row): "Shared" label repeated in D/E/F ---
$ws.Range("D8").Value = "Shared"
$ws.Range("E8").Value = "Shared"
$ws.Range("F8").Value = "Shared"
$ws.Range("D8:F8").Font.Bold = $true

# --- CPU rows 9-12: new columns D/E/F each = 1 ---
$ws.Range("D9").Value = 1
$ws.Range("E9").Value = 1
$ws.Range("F9").Value = 1

$ws.Range("D10").Value = 1
$ws.Range("E10").Value = 1
$ws.Range("F10").Value = 1

$ws.Range("D11").Value = 1
$ws.Range("E11").Value = 1
$ws.Range("F11").Value = 1

$ws.Range("D12").Value = 1
$ws.Range("E12").Value = 1
$ws.Range("F12").Value = 1

# --- Row 14: Max Slow formulas for D/E/F ---
$ws.Range("D14").Formula = "=MAX((D3/D9), (D4/D10), (D5/D11), (D6/D12))"
$ws.Range("E14").Formula = "=MAX((E3/E9), (E4/E10), (E5/E11), (E6/E12))"
$ws.Range("F14").Formula = "=MAX((F3/F9), (F4/F10), (F5/F11), (F6/F12))"

# --- Row 15: Wght Speed formulas for D/E/F ---
$ws.Range("D15").Formula = "=(D9/D3)+(D10/D4)+(D11/D5)+(D12/D6)"
$ws.Range("E15").Formula = "=(E9/E3)+(E10/E4)+(E11/E5)+(E12/E6)"
$ws.Range("F15").Formula = "=(F9/F3)+(F10/F4)+(F11/F5)+(F12/F6)"

# --- Row 17: NORM / Max Slow formulas for D/E/F ---
$ws.Range("D17").Formula = '=D14/$C$14'
$ws.Range("E17").Formula = '=E14/$C$14'
$ws.Range("F17").Formula = '=F14/$C$14'

# --- Row 18: Wght Speed formulas for D/E/F ---
$ws.Range("D18").Formula = '=D15/$C$15'
$ws.Range("E18").Formula = '=E15/$C$15'
$ws.Range("F18").Formula = '=F15/$C$15'

# --- Row 1: new blank formatted cell E1 (bold style) ---
$ws.Range("E1").Font.Bold = $true

# --- Update selection to match target state ---
$ws.Range("L12").Select()
